$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 110
$ws.Range("B110").Value = 6664683
$ws.Range("F110").Value = "Cerro Porteno"
$ws.Range("G110").Value = "General Caballero JLM"
$ws.Range("I110").Value = 2
$ws.Range("J110").Value = "A"
$ws.Range("K110").Value = 1.363
$ws.Range("L110").Value = 4.5
$ws.Range("M110").Value = 7
$ws.Range("N110").Value = 1.45
$ws.Range("O110").Value = 4.2
$ws.Range("P110").Value = 6
$ws.Range("Q110").Value = -1.25
$ws.Range("R110").Value = 2
$ws.Range("S110").Value = 1.8
$ws.Range("T110").Value = 2.75
$ws.Range("U110").Value = 1.775
$ws.Range("V110").Value = 2.025
$ws.Range("W110").Value = -1
$ws.Range("Y110").Value = 5
$ws.Range("Z110").Value = -1
$ws.Range("AA110").Value = 0.8
$ws.Range("AB110").Value = 0.3875
$ws.Range("AC110").Value = -0.5

# Row 111
$ws.Range("B111").Value = 6664682
$ws.Range("F111").Value = "Libertad Asuncion"
$ws.Range("G111").Value = "Olimpia Asuncion"
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = "H"
$ws.Range("K111").Value = 1.95
$ws.Range("L111").Value = 3.3
$ws.Range("M111").Value = 3.5
$ws.Range("N111").Value = 1.7
$ws.Range("O111").Value = 3.6
$ws.Range("P111").Value = 4.5
$ws.Range("Q111").Value = -0.75
$ws.Range("R111").Value = 1.9
$ws.Range("S111").Value = 1.9
$ws.Range("T111").Value = 2.5
$ws.Range("U111").Value = 1.85
$ws.Range("V111").Value = 1.95
$ws.Range("W111").Value = 0.7
$ws.Range("Y111").Value = -1
$ws.Range("Z111").Value = 0.45
$ws.Range("AA111").Value = -0.5
$ws.Range("AB111").Value = -1
$ws.Range("AC111").Value = 0.95

# Row 250
$ws.Range("B250").Value = 7454842
$ws.Range("F250").Value = "Sportivo Luqueno"
$ws.Range("G250").Value = "Libertad Asuncion"
$ws.Range("H250").Value = 0
$ws.Range("J250").Value = "A"
$ws.Range("K250").Value = 4
$ws.Range("L250").Value = 3.6
$ws.Range("M250").Value = 1.727
$ws.Range("N250").Value = 3.5
$ws.Range("O250").Value = 3.3
$ws.Range("P250").Value = 1.95
$ws.Range("Q250").Value = 0.5
$ws.Range("R250").Value = 1.8
$ws.Range("S250").Value = 2
$ws.Range("T250").Value = 2.5
$ws.Range("U250").Value = 1.975
$ws.Range("V250").Value = 1.825
$ws.Range("X250").Value = -1
$ws.Range("Y250").Value = 0.95
$ws.Range("AA250").Value = 1
$ws.Range("AC250").Value = 0.825

# Row 251
$ws.Range("B251").Value = 7453204
$ws.Range("F251").Value = "Cerro Porteno"
$ws.Range("G251").Value = "Tacuary"
$ws.Range("H251").Value = 1
$ws.Range("J251").Value = "D"
$ws.Range("K251").Value = 1.285
$ws.Range("L251").Value = 5
$ws.Range("M251").Value = 8
$ws.Range("N251").Value = 1.285
$ws.Range("O251").Value = 4.75
$ws.Range("P251").Value = 8
$ws.Range("Q251").Value = -1.5
$ws.Range("R251").Value = 1.9
$ws.Range("S251").Value = 1.9
$ws.Range("T251").Value = 3
$ws.Range("U251").Value = 1.9
$ws.Range("V251").Value = 1.9
$ws.Range("X251").Value = 3.75
$ws.Range("Y251").Value = -1
$ws.Range("AA251").Value = 0.8999999999999999
$ws.Range("AC251").Value = 0.8999999999999999

# Row 260
$ws.Range("B260").Value = 7493431
$ws.Range("F260").Value = "Sportivo Trinidense"
$ws.Range("G260").Value = "Guairena FC"
$ws.Range("H260").Value = 7
$ws.Range("J260").Value = "H"
$ws.Range("K260").Value = 2.05
$ws.Range("L260").Value = 3.3
$ws.Range("M260").Value = 3.3
$ws.Range("N260").Value = 2.6
$ws.Range("O260").Value = 3.1
$ws.Range("P260").Value = 2.6
$ws.Range("Q260").Value = 0
$ws.Range("R260").Value = 1.925
$ws.Range("S260").Value = 1.875
$ws.Range("T260").Value = 2.5
$ws.Range("U260").Value = 2
$ws.Range("V260").Value = 1.8
$ws.Range("W260").Value = 1.6
$ws.Range("Y260").Value = -1
$ws.Range("Z260").Value = 0.925
$ws.Range("AA260").Value = -1
$ws.Range("AB260").Value = 1
$ws.Range("AC260").Value = -1

# Row 261
$ws.Range("B261").Value = 7493310
$ws.Range("F261").Value = "Libertad Asuncion"
$ws.Range("G261").Value = "Tacuary"
$ws.Range("H261").Value = 1
$ws.Range("J261").Value = "A"
$ws.Range("K261").Value = 1.363
$ws.Range("L261").Value = 5
$ws.Range("M261").Value = 7
$ws.Range("N261").Value = 1.571
$ws.Range("O261").Value = 4.2
$ws.Range("P261").Value = 4.75
$ws.Range("Q261").Value = -0.75
$ws.Range("R261").Value = 1.8
$ws.Range("S261").Value = 2
$ws.Range("T261").Value = 2.75
$ws.Range("U261").Value = 1.8
$ws.Range("V261").Value = 2
$ws.Range("W261").Value = -1
$ws.Range("Y261").Value = 3.75
$ws.Range("Z261").Value = -1
$ws.Range("AA261").Value = 1
$ws.Range("AB261").Value = 0.4
$ws.Range("AC261").Value = -0.5

# Row 263
$ws.Range("B263").Value = 7493311
$ws.Range("F263").Value = "General Caballero JLM"
$ws.Range("G263").Value = "Olimpia Asuncion"
$ws.Range("H263").Value = 0
$ws.Range("J263").Value = "A"
$ws.Range("K263").Value = 3.4
$ws.Range("L263").Value = 3.3
$ws.Range("M263").Value = 2
$ws.Range("N263").Value = 3.2
$ws.Range("O263").Value = 3.25
$ws.Range("P263").Value = 2.1
$ws.Range("R263").Value = 1.95
$ws.Range("S263").Value = 1.85
$ws.Range("U263").Value = 1.775
$ws.Range("V263").Value = 2.025
$ws.Range("X263").Value = -1
$ws.Range("Y263").Value = 1.1
$ws.Range("Z263").Value = -1
$ws.Range("AA263").Value = 0.8500000000000001
$ws.Range("AB263").Value = -1
$ws.Range("AC263").Value = 1.025

# Row 264
$ws.Range("B264").Value = 7493433
$ws.Range("F264").Value = "Sportivo Luqueno"
$ws.Range("G264").Value = "Nacional Asuncion"
$ws.Range("H264").Value = 1
$ws.Range("I264").Value = 1
$ws.Range("J264").Value = "D"
$ws.Range("K264").Value = 2.75
$ws.Range("L264").Value = 3.2
$ws.Range("M264").Value = 2.4
$ws.Range("N264").Value = 2.75
$ws.Range("O264").Value = 3.1
$ws.Range("P264").Value = 2.45
$ws.Range("Q264").Value = 0.25
$ws.Range("R264").Value = 1.75
$ws.Range("S264").Value = 2.05
$ws.Range("T264").Value = 2.25
$ws.Range("U264").Value = 2
$ws.Range("V264").Value = 1.8
$ws.Range("W264").Value = -1
$ws.Range("X264").Value = 2.1
$ws.Range("Z264").Value = 0.375
$ws.Range("AA264").Value = -0.5
$ws.Range("AB264").Value = -0.5
$ws.Range("AC264").Value = 0.4

# Row 265
$ws.Range("B265").Value = 7493312
$ws.Range("F265").Value = "Cerro Porteno"
$ws.Range("G265").Value = "Guarani Asuncion"
$ws.Range("H265").Value = 4
$ws.Range("I265").Value = 0
$ws.Range("J265").Value = "H"
$ws.Range("K265").Value = 1.7
$ws.Range("L265").Value = 3.6
$ws.Range("M265").Value = 4.333
$ws.Range("N265").Value = 1.727
$ws.Range("O265").Value = 3.75
$ws.Range("P265").Value = 4.2
$ws.Range("Q265").Value = -0.5
$ws.Range("R265").Value = 1.8
$ws.Range("S265").Value = 2
$ws.Range("T265").Value = 2.75
$ws.Range("U265").Value = 1.875
$ws.Range("V265").Value = 1.925
$ws.Range("W265").Value = 0.7270000000000001
$ws.Range("Y265").Value = -1
$ws.Range("Z265").Value = 0.8
$ws.Range("AA265").Value = -1
$ws.Range("AB265").Value = 0.875
$ws.Range("AC265").Value = -1

# Row 301
$ws.Range("O301").Value = 3.6
$ws.Range("P301").Value = 4.333
$ws.Range("T301").Value = 2.5
$ws.Range("U301").Value = 1.975
$ws.Range("V301").Value = 1.825

# Row 302
$ws.Range("N302").Value = 2.1
$ws.Range("O302").Value = 3.3
$ws.Range("P302").Value = 3.1
$ws.Range("R302").Value = 1.85
$ws.Range("S302").Value = 1.95
$ws.Range("U302").Value = 1.825
$ws.Range("V302").Value = 1.975

# Row 304
$ws.Range("R304").Value = 1.975
$ws.Range("S304").Value = 1.825

# Row 305
$ws.Range("R305").Value = 2.075
$ws.Range("S305").Value = 1.725
